# Fruta / hortaliza, semanal
#
# A new daily-price record is inserted as row 62 (Vega Central Mapocho de
# Santiago, Frambuesa, Primera, Provincia de Curicó) which pushes every
# existing record from the old row 62 down through the old row 118 down by
# one row (old row 118 becomes row 119). The sheet's used range grows from
# A1:T118 to A1:T119.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 62, shifting rows 62:118 down to 63:119.
$ws.Rows(62).Insert()

# Populate the newly inserted row 62 with the new record.
$ws.Cells.Item(62, 1).Value = 9
$ws.Cells.Item(62, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = 44944
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100101
$ws.Cells.Item(62, 8).Value = "Berries"
$ws.Cells.Item(62, 9).Value = 100101004
$ws.Cells.Item(62, 10).Value = "Frambuesa"
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 280
$ws.Cells.Item(62, 14).Value = 8000
$ws.Cells.Item(62, 15).Value = 8000
$ws.Cells.Item(62, 16).Value = 8000
$ws.Cells.Item(62, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(62, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(62, 19).Value = 4000
$ws.Cells.Item(62, 20).Value = 2
